# Add NSW "Anzac Day (additional)" observance rows for 2026 and 2027.
#
# Anzac Day (25 Apr) falls on Saturday in 2026 (observed Mon 27 Apr) and on
# Sunday in 2027 (observed Mon 26 Apr). NSW adds an extra "Anzac Day
# (additional)" holiday on the observed Monday in those cases. Insert a new
# row directly after each NSW "Anzac Day" row (pushing subsequent rows down),
# and populate it with the observed date / state / name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2026: Anzac Day (Sat 25 Apr) -> additional observed Mon 27 Apr ---
# NSW "Anzac Day" currently sits at row 15; insert the new row at row 16
# (pushing the VIC "Anzac Day" row, and everything after, down by one).
# The date column holds plain text (e.g. "2026-04-25"), so the value is
# entered with a leading apostrophe to stop Excel from auto-converting the
# date-shaped text into a date serial number, matching the existing cells.
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = "'2026-04-27"
$ws.Range("B16").Value = "NSW"
$ws.Range("C16").Value = "Anzac Day (additional)"

# --- 2027: Anzac Day (Sun 25 Apr) -> additional observed Mon 26 Apr ---
# After the insertion above, NSW "Anzac Day" 2027 now sits at row 43; insert
# the new row at row 44 (pushing the VIC "Anzac Day" row, and everything
# after, down by one).
$ws.Rows.Item(44).Insert()
$ws.Range("A44").Value = "'2027-04-26"
$ws.Range("B44").Value = "NSW"
$ws.Range("C44").Value = "Anzac Day (additional)"
